$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2545.2
$ws.Range("I86").Value = 2418.5454
$ws.Range("J86").Value = 2700
$ws.Range("K86").Value = 2418.5454
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -1295.5454
$ws.Range("N86").Value = -4946

$ws.Range("H89").Value = 2545.2
$ws.Range("I89").Value = 2418.5454
$ws.Range("J89").Value = 2700
$ws.Range("K89").Value = 12092.727
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -6476.726999999999
$ws.Range("N89").Value = -24732

$ws.Range("H112").Value = 1947.4902
$ws.Range("J112").Value = 2072.1956
$ws.Range("L112").Value = 6216.5868
$ws.Range("N112").Value = -8432.586800000001

$ws.Range("H129").Value = 991.9655
$ws.Range("I129").Value = 471.5
$ws.Range("K129").Value = 1414.5
$ws.Range("M129").Value = 3585.5

$ws.Range("H138").Value = 2115.74
$ws.Range("I138").Value = 1209.3226
$ws.Range("J138").Value = 2522.971
$ws.Range("K138").Value = 3627.9678
$ws.Range("L138").Value = 7568.913
$ws.Range("M138").Value = 1512.0322
$ws.Range("N138").Value = -17848.913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16981.363
$ws.Range("I32").Value = 13440.661
$ws.Range("J32").Value = 48341.855
$ws.Range("K32").Value = 13440.661
$ws.Range("L32").Value = 48341.855
$ws.Range("M32").Value = -13153.661
$ws.Range("N32").Value = -48915.855

$ws.Range("H97").Value = 638.125
$ws.Range("I97").Value = 638.125
$ws.Range("K97").Value = 638.125
$ws.Range("M97").Value = -142.125

$ws.Range("H101").Value = 56602
$ws.Range("J101").Value = 56602
$ws.Range("L101").Value = 56602
$ws.Range("N101").Value = -63092

$ws.Range("H102").Value = 2010
$ws.Range("I102").Value = 2010
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2010
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -388

$ws.Range("H105").Value = 500370
$ws.Range("J105").Value = 500370
$ws.Range("L105").Value = 500370
$ws.Range("N105").Value = -507358

$ws.Range("H106").Value = 44900
$ws.Range("J106").Value = 44900
$ws.Range("L106").Value = 44900
$ws.Range("N106").Value = -47424

$ws.Range("H132").Value = 1675955
$ws.Range("I132").Value = 2002046.1
$ws.Range("J132").Value = 45499.5
$ws.Range("K132").Value = 6006138.300000001
$ws.Range("L132").Value = 136498.5
$ws.Range("M132").Value = -6003608.300000001
$ws.Range("N132").Value = -141558.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 261882.5
$ws.Range("I92").Value = 23364
$ws.Range("J92").Value = 500401
$ws.Range("K92").Value = 23364
$ws.Range("L92").Value = 500401
$ws.Range("M92").Value = -20868
$ws.Range("N92").Value = -505393

$ws.Range("H103").Value = 21900
$ws.Range("J103").Value = 21900
$ws.Range("L103").Value = 21900
$ws.Range("N103").Value = -24244

$ws.Range("H134").Value = 419333.78
$ws.Range("I134").Value = 608109.2
$ws.Range("J134").Value = 4027.8667
$ws.Range("K134").Value = 1824327.6
$ws.Range("L134").Value = 12083.6001
$ws.Range("M134").Value = -1821792.6
$ws.Range("N134").Value = -17153.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3705.4614
$ws.Range("I31").Value = 2638.1177
$ws.Range("J31").Value = 5721.5557
$ws.Range("K31").Value = 2638.1177
$ws.Range("L31").Value = 5721.5557
$ws.Range("M31").Value = -2343.1177
$ws.Range("N31").Value = -6311.5557

$ws.Range("H34").Value = 3705.4614
$ws.Range("I34").Value = 2638.1177
$ws.Range("J34").Value = 5721.5557
$ws.Range("K34").Value = 2638.1177
$ws.Range("L34").Value = 5721.5557
$ws.Range("M34").Value = -2436.1177
$ws.Range("N34").Value = -6125.5557

$ws.Range("H58").Value = 1755.0454
$ws.Range("I58").Value = 1647.8889
$ws.Range("J58").Value = 2237.25
$ws.Range("K58").Value = 1647.8889
$ws.Range("L58").Value = 2237.25
$ws.Range("M58").Value = -1444.8889
$ws.Range("N58").Value = -2643.25

$ws.Range("H92").Value = 64866.332
$ws.Range("J92").Value = 64866.332
$ws.Range("L92").Value = 64866.332
$ws.Range("N92").Value = -69858.33199999999

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0

$ws.Range("H106").Value = 500000
$ws.Range("J106").Value = 500000
$ws.Range("L106").Value = 500000
$ws.Range("N106").Value = -502524

$ws.Range("H109").Value = 19997.5
$ws.Range("J109").Value = 19997.5
$ws.Range("L109").Value = 19997.5
$ws.Range("N109").Value = -22077.5

$ws.Range("H136").Value = 1755.0454
$ws.Range("I136").Value = 1647.8889
$ws.Range("J136").Value = 2237.25
$ws.Range("K136").Value = 4943.6667
$ws.Range("L136").Value = 6711.75
$ws.Range("M136").Value = -2393.6667
$ws.Range("N136").Value = -11811.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3994.8572
$ws.Range("J5").Value = 1405
$ws.Range("L5").Value = 4215
$ws.Range("N5").Value = -4439

$ws.Range("H108").Value = 1519
$ws.Range("I108").Value = 1013.5
$ws.Range("J108").Value = 2530
$ws.Range("K108").Value = 3040.5
$ws.Range("L108").Value = 7590
$ws.Range("M108").Value = -160.5
$ws.Range("N108").Value = -13350

$ws.Range("H109").Value = 1685.2727
$ws.Range("I109").Value = 663.8333
$ws.Range("J109").Value = 2911
$ws.Range("K109").Value = 1991.4999
$ws.Range("L109").Value = 8733
$ws.Range("M109").Value = -951.4999
$ws.Range("N109").Value = -10813

$ws.Range("H112").Value = 3628.6428
$ws.Range("J112").Value = 4394.3184
$ws.Range("L112").Value = 13182.9552
$ws.Range("N112").Value = -15398.9552

$ws.Range("H115").Value = 2567.1667
$ws.Range("I115").Value = 339.33334
$ws.Range("J115").Value = 3012.7334
$ws.Range("K115").Value = 1018.00002
$ws.Range("L115").Value = 9038.200199999999
$ws.Range("M115").Value = 156.9999799999999
$ws.Range("N115").Value = -11388.2002

$ws.Range("H118").Value = 1879.1538
$ws.Range("I118").Value = 1104.1428
$ws.Range("J118").Value = 2783.3333
$ws.Range("K118").Value = 3312.4284
$ws.Range("L118").Value = 8349.999899999999
$ws.Range("M118").Value = -2069.4284
$ws.Range("N118").Value = -10835.9999

$ws.Range("H121").Value = 1350.5385
$ws.Range("I121").Value = 720.6667
$ws.Range("J121").Value = 1890.4286
$ws.Range("K121").Value = 2162.0001
$ws.Range("L121").Value = 5671.2858
$ws.Range("M121").Value = -852.0001000000002
$ws.Range("N121").Value = -8291.2858

$ws.Range("H122").Value = 814.46155
$ws.Range("I122").Value = 421
$ws.Range("J122").Value = 1273.5
$ws.Range("K122").Value = 3789
$ws.Range("L122").Value = 11461.5
$ws.Range("M122").Value = -1339
$ws.Range("N122").Value = -16361.5

$ws.Range("H123").Value = 7831.1763
$ws.Range("J123").Value = 8733.333000000001
$ws.Range("L123").Value = 26199.999
$ws.Range("N123").Value = -31099.999

$ws.Range("H125").Value = 3208
$ws.Range("J125").Value = 3993.3333
$ws.Range("L125").Value = 11979.9999
$ws.Range("N125").Value = -21819.9999

$ws.Range("H132").Value = 1960.8
$ws.Range("I132").Value = 1143.1428
$ws.Range("J132").Value = 2676.25
$ws.Range("K132").Value = 10288.2852
$ws.Range("L132").Value = 24086.25
$ws.Range("M132").Value = -7758.2852
$ws.Range("N132").Value = -29146.25

$ws.Range("H133").Value = 3153.7693
$ws.Range("I133").Value = 1811.8
$ws.Range("J133").Value = 3992.5
$ws.Range("K133").Value = 5435.4
$ws.Range("L133").Value = 11977.5
$ws.Range("M133").Value = -375.3999999999996
$ws.Range("N133").Value = -22097.5

$ws.Range("H134").Value = 4419.7715
$ws.Range("I134").Value = 2374.6843
$ws.Range("J134").Value = 6848.3125
$ws.Range("K134").Value = 7124.0529
$ws.Range("L134").Value = 20544.9375
$ws.Range("M134").Value = -2054.0529
$ws.Range("N134").Value = -30684.9375

$ws.Range("H135").Value = 3994.8572
$ws.Range("J135").Value = 1405
$ws.Range("L135").Value = 12645
$ws.Range("N135").Value = -17715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 73941.42999999999
$ws.Range("I97").Value = 73941.42999999999
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 73941.42999999999
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -73445.42999999999

$ws.Range("H123").Value = 11204.546
$ws.Range("J123").Value = 11204.546
$ws.Range("L123").Value = 11204.546
$ws.Range("N123").Value = -16104.546

$ws.Range("H132").Value = 2612.6365
$ws.Range("I132").Value = 1842.5
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 5527.5
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -2997.5
$ws.Range("N132").Value = -19059.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5333.3335
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9888

$ws.Range("H126").Value = 5333.3335
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 8666.5
$ws.Range("J6").Value = 8400
$ws.Range("L6").Value = 8400
$ws.Range("N6").Value = -8630

$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15500

$ws.Range("H36").Value = 15000
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15500

$ws.Range("H39").Value = 9857235
$ws.Range("I39").Value = 19699470
$ws.Range("K39").Value = 19699470
$ws.Range("M39").Value = -19699057

$ws.Range("H48").Value = 9029.5
$ws.Range("I48").Value = 3059
$ws.Range("K48").Value = 3059
$ws.Range("M48").Value = -2490

$ws.Range("H132").Value = 1955.3243
$ws.Range("I132").Value = 1300.9656
$ws.Range("J132").Value = 4327.375
$ws.Range("K132").Value = 3902.8968
$ws.Range("L132").Value = 12982.125
$ws.Range("M132").Value = -1372.8968
$ws.Range("N132").Value = -18042.125
